$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..40 down to 4..41.
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the same metadata as row 2 (same market/product),
# but with the new date / price observation.
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44882
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100108
$ws.Range("H3").Value = "Tropicales y subtropicales"
$ws.Range("I3").Value = 100108007
$ws.Range("J3").Value = "Coco"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 28000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29000
$ws.Range("Q3").Value = '$/malla 20 unidades'
$ws.Range("R3").Value = "Perú"
$ws.Range("S3").Value = 1450
$ws.Range("T3").Value = 20
